# Correcting Relevance Markers Appenzeller-Herzog (2019) - van Dis (2020)
#
# Updates the "record_atd" (time-to-discovery) simulation results for the
# second simulation run (td_sim_1, column C) and the resulting average
# (record_atd, column D) for rows 2-112, plus the overall average of
# td_sim_1 in row 113.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column C (td_sim_1) for rows 2..112
$newC = @(
    1891,1320,3682,4389,3193,1570,1459,1641,1619,3142,
    2900,501,1014,3469,1034,1617,1292,2003,3723,982,
    1938,2552,736,3183,4549,3563,1523,2002,3637,1895,
    1936,1696,934,1218,3850,1373,1536,130,1132,3426,
    392,2484,984,773,823,1231,531,5126,4588,1529,
    2197,2386,2976,2042,482,774,2858,1943,2909,1975,
    2021,535,466,1027,4967,905,800,2698,2370,3279,
    1656,1398,454,129,538,2090,1899,1507,476,785,
    2029,2454,3724,114,2226,4233,2958,811,1085,572,
    2999,2395,1210,3153,3137,3797,1628,977,1662,1687,
    1535,2771,1453,3933,2353,670,4570,2018,2152,2778,
    3172
)

# Corresponding recomputed values for column D (record_atd = AVERAGE(B,C))
$newD = @(
    1137,1325,3661.5,4276,3146.5,1573.5,1700.5,1571,1595.5,3045,
    2933,360,1046,3462,1065,1593.5,1298.5,1885.5,3650.5,999,
    1196.5,2513.5,812,3115,4494.5,3496.5,1484,1998.5,3670.5,1164,
    1628.5,1587,900,1244.5,3891,1371,1524.5,460.5,1039,3413.5,
    320.5,2355,1000.5,644,711.5,1250,341.5,5125.5,4557.5,1539.5,
    2205.5,1523.5,2853.5,1943.5,297.5,645,2852,1207,3022,1845.5,
    2018,345.5,466,1063.5,4970.5,794,669,2701.5,2409.5,3231,
    1656.5,1409.5,241.5,462,348,2066,1182,1535,291.5,676,
    2013,2353,3652,481.5,2185.5,4230,2974,678,1105,381.5,
    2828.5,2385.5,1224.5,3080,3027.5,3841,1598.5,999,1650.5,1665,
    1535,2782,1660,3748,2399.5,496,4538.5,1898,2131,2749.5,
    3090.5
)

$startRow = 2
for ($i = 0; $i -lt $newC.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 3).Value = $newC[$i]
    $ws.Cells.Item($row, 4).Value = $newD[$i]
}

# Overall average of td_sim_1 (row 113) recomputed from the new column C values
$ws.Cells.Item(113, 3).Value = 2022.603603603604
